# Update "想去人数" (want-to-go count) figures across the workbook's
# sheets to reflect the latest scrape (commit: "Update gh-pages to output
# generated at 456a3b4").
#
# Column F on each of the data sheets holds this count. Sheet "本地生活"
# is not touched by this update.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 9024
$ws.Range("F3").Value = 2699
$ws.Range("F5").Value = 842
$ws.Range("F6").Value = 740
$ws.Range("F7").Value = 139
$ws.Range("F8").Value = 82
$ws.Range("F9").Value = 389
$ws.Range("F10").Value = 918
$ws.Range("F11").Value = 4000
$ws.Range("F12").Value = 316
$ws.Range("F13").Value = 197
$ws.Range("F14").Value = 814
$ws.Range("F15").Value = 781
$ws.Range("F17").Value = 510
$ws.Range("F20").Value = 1447
$ws.Range("F21").Value = 1371
$ws.Range("F22").Value = 520
$ws.Range("F24").Value = 159
$ws.Range("F25").Value = 185
$ws.Range("F26").Value = 385
$ws.Range("F27").Value = 77
$ws.Range("F28").Value = 1024
$ws.Range("F31").Value = 772
$ws.Range("F32").Value = 77
$ws.Range("F34").Value = 112
$ws.Range("F39").Value = 207
$ws.Range("F40").Value = 428
$ws.Range("F41").Value = 35

# --- Sheet: 演出 ---------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 26
$ws.Range("F4").Value = 79
$ws.Range("F6").Value = 55
$ws.Range("F7").Value = 3

# --- Sheet: 全部类型 ------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 9024
$ws.Range("F5").Value = 842
$ws.Range("F6").Value = 740
$ws.Range("F7").Value = 139
$ws.Range("F8").Value = 82
$ws.Range("F9").Value = 389
$ws.Range("F10").Value = 918
$ws.Range("F11").Value = 26
$ws.Range("F12").Value = 4000
$ws.Range("F13").Value = 316
$ws.Range("F14").Value = 197
$ws.Range("F16").Value = 79
$ws.Range("F17").Value = 814
$ws.Range("F18").Value = 781
$ws.Range("F20").Value = 55
$ws.Range("F22").Value = 510
$ws.Range("F24").Value = 3
$ws.Range("F26").Value = 1447
$ws.Range("F27").Value = 1371
$ws.Range("F28").Value = 520
$ws.Range("F30").Value = 159
$ws.Range("F31").Value = 185
$ws.Range("F33").Value = 385
$ws.Range("F34").Value = 77
$ws.Range("F35").Value = 1024
$ws.Range("F37").Value = 772
$ws.Range("F38").Value = 77
$ws.Range("F40").Value = 112
$ws.Range("F44").Value = 207
$ws.Range("F45").Value = 428
$ws.Range("F46").Value = 35
